$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Modelo" header in F1, reusing the exact same header
# formatting (bold font, border, centered) that E1 already has.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Modelo"

# Update the computed metric values in row 2 (re-run produced slightly
# different MSE / R2 / MAE)
$ws.Range("B2").Value = 0.5389788899637855
$ws.Range("C2").Value = 0.9892677857137916
$ws.Range("D2").Value = 0.6112735950122063

# Add the new model name value in F2
$ws.Range("F2").Value = "Pipeline(steps=[('model', RandomForestRegressor(max_depth=5))])"
